$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scores")

$ws.Range("C2").Value = "Fairly Glib"
$ws.Range("C3").Value = "Fairly Grandiose"
$ws.Range("C4").Value = "Fairly Conniving"
$ws.Range("C5").Value = "Fairly Deceptive"
$ws.Range("C6").Value = "Fairly Unremorseful"
$ws.Range("C7").Value = "Fairly Callous"
$ws.Range("C8").Value = "Fairly Inexpressive"
$ws.Range("C9").Value = "Fairly Irresponsible"
$ws.Range("C10").Value = "Fairly Sensation Seeking"
$ws.Range("C11").Value = "Fairly Unrealistic"
$ws.Range("C12").Value = "Fairly Impulsive"
$ws.Range("C13").Value = "Fairly Irresponsible"
$ws.Range("C14").Value = "Fairly Parasitic"
$ws.Range("C15").Value = "Fairly Noncommittal"
$ws.Range("C16").Value = "Fairly Promiscuous"
$ws.Range("C17").Value = "Fairly Emotionally Controlled"
$ws.Range("C18").Value = "Fairly Problematic"
$ws.Range("C19").Value = "Fairly Delinquent"
$ws.Range("C20").Value = "Fairly Noncompliant"
$ws.Range("C21").Value = "Fairly Versatile"
